# Auto-generated update of cryptos list values (price / volume / ranking)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.351.83"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.719.32"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.14"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5313"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06706"
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2663"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.96"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07707"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.483"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.954.41"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.713.80"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5832"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8218"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.23"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.393.33"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.14"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.664"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.48"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.046"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.708"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1209"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.253"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.29"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05388"
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.297"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.481"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.431"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.638"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.863"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9531"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.399"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5895"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.155.29"
$ws.Range("E39").Value = "  +10.12%  "
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.807"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.007"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8416"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.861.08"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("E46").Value = "  -6.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.90"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.162"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05205"
$ws.Range("E51").Value = "  -0.83%  "
